# users_db.xlsx — add a "gender" column (derived from first_name) right
# before "city", move the selection, and nudge the saved window position.
#
# Column layout BEFORE: rut, first_name, last_name, email, cellphone,
#   city, card_status, transdata_id, place                (A..I)
# Column layout AFTER:  rut, first_name, last_name, email, cellphone,
#   gender, city, card_status, transdata_id, place         (A..J)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at F; everything from the old F ("city") onward
# shifts one column to the right automatically (F->G, G->H, H->I, I->J).
$ws.Columns.Item(6).Insert()

# New header
$ws.Cells.Item(1, 6).Value = "gender"

# Per-row gender values (row number -> gender), matching first_name in
# column B for each record.
$genderByRow = @{
  2  = "male";   3  = "male";   4  = "male";   5  = "male";
  6  = "male";   7  = "male";   8  = "male";   9  = "female";
  10 = "female"; 11 = "female"; 12 = "female"; 13 = "male";
  14 = "female"; 15 = "female"; 16 = "male";   17 = "male";
  18 = "male";   19 = "male";   20 = "female"; 21 = "male";
  22 = "male";   23 = "female"; 24 = "male";   25 = "female";
  26 = "male";   27 = "female"; 28 = "male";   29 = "male";
  30 = "male";   31 = "male";   32 = "male";   33 = "female";
  34 = "male";   35 = "male";   36 = "male";   37 = "male";
  38 = "male"
}

foreach ($row in 2..38) {
  $ws.Cells.Item($row, 6).Value = $genderByRow[$row]
}

# Move the selection / scroll position to match the edited workbook.
$ws.Range("F38").Select()
try { $excel.ActiveWindow.ScrollRow = 11 } catch {}
try { $excel.ActiveWindow.ScrollColumn = 1 } catch {}

# Best-effort: nudge the saved window position (xWindow/yWindow).
try {
  $win = $wb.Windows.Item(1)
  $win.Left = 340
  $win.Top = 480
} catch {}
try {
  $excel.ActiveWindow.Left = 340
  $excel.ActiveWindow.Top = 480
} catch {}
